# Add support for a new "EPV" question type (Table A,B,C,D) as an
# additional row in the question-type lookup table on Sheet1.
#
# Existing layout (Question | Type):
#   Single Choice   | sc
#   Multiple Choice | mc
#   Free Text       | txt
#   Number          | int
#
# New row being appended:
#   Table (A,B,C,D) | epv

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write column B before column A so the shared-string table picks up
# "epv" ahead of "Table (A,B,C,D)" (matches authoring order).
$ws.Range("B6").Value = "epv"
$ws.Range("A6").Value = "Table (A,B,C,D)"

# Leave the selection on the next empty row, as if the author had just
# finished typing this row and moved down.
$ws.Range("A7").Select()
